$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

$styled = $excel.Union($ws.Range("B1"), $ws.Range("A2"))
$styled.Font.Bold = $true
$styled.HorizontalAlignment = -4108
$styled.VerticalAlignment = -4160
$styled.Borders.LineStyle = 1
